# LLTD_sweep.xlsx edit script
# Summary of the change (per commit message: "Fixed the problem of peak
# longitudinal acceleration in the start / Tried running Budapest"):
#   - MassInertia: relabel the 3 alternate setups from "-5/-10/-15% lltd"
#     to "3/6/9 lltd"
#   - Aero: front downforce coefficient (col B) 3.5 -> 4.5
#   - Susp: LLTD sweep values changed (0.6/0.55/0.5/0.45 -> 0.5/0.53/0.56/0.59)
#   - Brake: max torque bumped from ~5000 to 8000 (flat across all 4 setups)
#   - Engine: gear ratio changed (3 -> 3.5), power drag bumped 450 -> 550,
#     and a new "gear_ratio" column pointing at GearRatio.xlsx was added
#   - Aero sheet left as the active/selected tab

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# MassInertia: rename the sweep setups in column B (rows 3-5)
# ---------------------------------------------------------------------
$wsMass = $wb.Worksheets.Item("MassInertia")
$wsMass.Range("B3").Value = "3 lltd"
$wsMass.Range("B4").Value = "6 lltd"
$wsMass.Range("B5").Value = "9 lltd"

# ---------------------------------------------------------------------
# Aero: bump column B (aerobalance-related input) from 3.5 to 4.5
# ---------------------------------------------------------------------
$wsAero = $wb.Worksheets.Item("Aero")
$wsAero.Range("B2:B5").Value = 4.5

# ---------------------------------------------------------------------
# Susp: rewrite the LLTD sweep column
# ---------------------------------------------------------------------
$wsSusp = $wb.Worksheets.Item("Susp")
$wsSusp.Range("A2").Value = 0.5
$wsSusp.Range("A3").Value = 0.53
$wsSusp.Range("A4").Value = 0.56
$wsSusp.Range("A5").Value = 0.59

# ---------------------------------------------------------------------
# Brake: all four setups now share the same max torque, 8000
# ---------------------------------------------------------------------
$wsBrake = $wb.Worksheets.Item("Brake")
$wsBrake.Range("B2:B5").Value = 8000

# ---------------------------------------------------------------------
# Engine: gear ratio 3 -> 3.5, power-drag 450 -> 550, plus a new
# gear_ratio column referencing GearRatio.xlsx
# ---------------------------------------------------------------------
$wsEngine = $wb.Worksheets.Item("Engine")
$wsEngine.Range("A2:A5").Value = 3.5
$wsEngine.Range("C2:C5").Value = 550
$wsEngine.Range("E1").Value = "gear_ratio"
$wsEngine.Range("E2:E5").Value = "GearRatio.xlsx"
$wsEngine.Columns.Item(5).ColumnWidth = 13.44140625
[void]$wsEngine.Range("A4:E4").Select()

# ---------------------------------------------------------------------
# Selection / active-tab bookkeeping to mirror the saved UI state
# ---------------------------------------------------------------------
[void]$wsMass.Range("A6").Select()
[void]$wsBrake.Range("A5:B5").Select()

$wsAero.Activate()
[void]$wsAero.Range("G13").Select()
